$wb = $excel.ActiveWorkbook

# Sheets affected: "展览" and "全部类型" (both contain the same data rows)
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1653
    $ws.Range("F6").Value = 44
}
